# "Add files via upload" - re-upload of the Saldo export with a handful of
# account rows removed from the "Export" sheet.
#
# Rows removed (1-based worksheet row numbers, header is row 1):
#   4  -> 005105172 VALDIVINO   102848.76
#   8  -> 004363260 LARISSA      15324.44
#   9  -> 005242683 LUCAS        13428.59
#   11 -> 005064129 THIAGO        3187.11
#   12 -> 004482102 NATALIA       2000
#   13 -> 004363250 HELIO         1624.55
#   15 -> 002687737 JOSE          1418.72
#   16 -> 001759765 NATAL         1401.18

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsToDelete = @(4, 8, 9, 11, 12, 13, 15, 16) | Sort-Object -Descending

foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}
